# "update data with resort sheetname"
# Re-sort the worksheet tabs: "总计" moves in front of "2020-Q4" so that
# "总计" becomes the first sheet (was second) and "2020-Q4" becomes the
# second sheet (was first). No cell data changes - only the sheet order.

$wb = $excel.ActiveWorkbook

$sheetTotal = $wb.Worksheets.Item("总计")
$sheetQ4 = $wb.Worksheets.Item("2020-Q4")

# Move "总计" so it sits immediately before "2020-Q4" -> new order:
# [总计, 2020-Q4]
$sheetTotal.Move($sheetQ4)

# "2020-Q4" keeps being the selected/active tab (it was the active tab
# before the reorder), so re-select it by name now that positions changed.
$wb.Worksheets.Item("2020-Q4").Activate()
